$d = $word.ActiveDocument

# The doc has two floating (anchored) pictures living in the primary
# header: one PNG (currently named "image2.png") and one JPG (currently
# named "image1.jpg"). The commit swaps the cosmetic name <-> extension
# pairing: the PNG becomes "image1.png" and the JPG becomes "image2.jpg"
# (the underlying embedded media / relationships are untouched - only
# the shape's display name changes).

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    # wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3
    foreach ($idx in 1, 2, 3) {
        $hdr = $section.Headers.Item($idx)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Shapes.Count; $i++) {
                $shp = $hdr.Shapes.Item($i)
                if ($shp.Name -eq "image2.png") {
                    $shp.Name = "image1.png"
                } elseif ($shp.Name -eq "image1.jpg") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    foreach ($idx in 1, 2, 3) {
        $ftr = $section.Footers.Item($idx)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Shapes.Count; $i++) {
                $shp = $ftr.Shapes.Item($i)
                if ($shp.Name -eq "image2.png") {
                    $shp.Name = "image1.png"
                } elseif ($shp.Name -eq "image1.jpg") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
